$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3426.375
$ws.Range("I64").Value = 3089.1667
$ws.Range("J64").Value = 3538.7778
$ws.Range("K64").Value = 3089.1667
$ws.Range("L64").Value = 3538.7778
$ws.Range("M64").Value = -2841.1667
$ws.Range("N64").Value = -4034.7778

$ws.Range("H67").Value = 3426.375
$ws.Range("I67").Value = 3089.1667
$ws.Range("J67").Value = 3538.7778
$ws.Range("K67").Value = 3089.1667
$ws.Range("L67").Value = 3538.7778
$ws.Range("M67").Value = -2231.1667
$ws.Range("N67").Value = -5254.7778

$ws.Range("H74").Value = 5320.6924
$ws.Range("I74").Value = 5887.769
$ws.Range("J74").Value = 4753.615
$ws.Range("K74").Value = 5887.769
$ws.Range("L74").Value = 4753.615
$ws.Range("M74").Value = -4951.769
$ws.Range("N74").Value = -6625.615

$ws.Range("H76").Value = 3332.8572
$ws.Range("I76").Value = 3245
$ws.Range("J76").Value = 3450
$ws.Range("K76").Value = 3245
$ws.Range("L76").Value = 3450
$ws.Range("M76").Value = -2930
$ws.Range("N76").Value = -4080

$ws.Range("H77").Value = 5320.6924
$ws.Range("I77").Value = 5887.769
$ws.Range("J77").Value = 4753.615
$ws.Range("K77").Value = 29438.845
$ws.Range("L77").Value = 23768.075
$ws.Range("M77").Value = -24758.845
$ws.Range("N77").Value = -33128.075

$ws.Range("H79").Value = 3332.8572
$ws.Range("I79").Value = 3245
$ws.Range("J79").Value = 3450
$ws.Range("K79").Value = 3245
$ws.Range("L79").Value = 3450
$ws.Range("M79").Value = -2153
$ws.Range("N79").Value = -5634

$ws.Range("H127").Value = 1865.7587
$ws.Range("I127").Value = 1502.25
$ws.Range("J127").Value = 2004.238
$ws.Range("K127").Value = 4506.75
$ws.Range("L127").Value = 6012.714
$ws.Range("M127").Value = 453.25
$ws.Range("N127").Value = -15932.714

$ws.Range("H137").Value = 1092
$ws.Range("I137").Value = 924
$ws.Range("J137").Value = 1344
$ws.Range("K137").Value = 2772
$ws.Range("L137").Value = 4032
$ws.Range("M137").Value = -222
$ws.Range("N137").Value = -9132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2219.6875
$ws.Range("I63").Value = 2126.25
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2126.25
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1440.25
$ws.Range("N63").Value = -3872

$ws.Range("H66").Value = 2219.6875
$ws.Range("I66").Value = 2126.25
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 10631.25
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -7199.25
$ws.Range("N66").Value = -19364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2002.8837
$ws.Range("I31").Value = 1426.4073
$ws.Range("J31").Value = 2975.6875
$ws.Range("K31").Value = 1426.4073
$ws.Range("L31").Value = 2975.6875
$ws.Range("M31").Value = -1131.4073
$ws.Range("N31").Value = -3565.6875

$ws.Range("H34").Value = 2002.8837
$ws.Range("I34").Value = 1426.4073
$ws.Range("J34").Value = 2975.6875
$ws.Range("K34").Value = 1426.4073
$ws.Range("L34").Value = 2975.6875
$ws.Range("M34").Value = -1224.4073
$ws.Range("N34").Value = -3379.6875

$ws.Range("H62").Value = 37040108
$ws.Range("I62").Value = 66669070
$ws.Range("J62").Value = 3901.5
$ws.Range("K62").Value = 66669070
$ws.Range("L62").Value = 3901.5
$ws.Range("M62").Value = -66668446
$ws.Range("N62").Value = -5149.5

$ws.Range("H65").Value = 37040108
$ws.Range("I65").Value = 66669070
$ws.Range("J65").Value = 3901.5
$ws.Range("K65").Value = 333345350
$ws.Range("L65").Value = 19507.5
$ws.Range("M65").Value = -333342230
$ws.Range("N65").Value = -25747.5

$ws.Range("H99").Value = 1892.3846
$ws.Range("I99").Value = 1727.3636
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1727.3636
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -229.3635999999999
$ws.Range("N99").Value = -5796

$ws.Range("H126").Value = 1892.3846
$ws.Range("I126").Value = 1727.3636
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 5182.0908
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -2712.0908
$ws.Range("N126").Value = -13340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 441.05884
$ws.Range("I5").Value = 393.16666
$ws.Range("J5").Value = 556
$ws.Range("K5").Value = 1179.49998
$ws.Range("L5").Value = 1668
$ws.Range("M5").Value = -1067.49998
$ws.Range("N5").Value = -1892

$ws.Range("H131").Value = 947.747
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 1011.8219
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 3035.4657
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -13115.4657

$ws.Range("H135").Value = 441.05884
$ws.Range("I135").Value = 393.16666
$ws.Range("J135").Value = 556
$ws.Range("K135").Value = 3538.49994
$ws.Range("L135").Value = 5004
$ws.Range("M135").Value = -1003.49994
$ws.Range("N135").Value = -10074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5081.657
$ws.Range("I70").Value = 4697.619
$ws.Range("J70").Value = 5657.7144
$ws.Range("K70").Value = 4697.619
$ws.Range("L70").Value = 5657.7144
$ws.Range("M70").Value = -4427.619
$ws.Range("N70").Value = -6197.7144

$ws.Range("H73").Value = 5081.657
$ws.Range("I73").Value = 4697.619
$ws.Range("J73").Value = 5657.7144
$ws.Range("K73").Value = 4697.619
$ws.Range("L73").Value = 5657.7144
$ws.Range("M73").Value = -3761.619
$ws.Range("N73").Value = -7529.7144

$ws.Range("H80").Value = 2616.6667
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2640
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2640
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4636

$ws.Range("H83").Value = 2616.6667
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2640
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 13200
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -23184

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 9000
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2520.5334
$ws.Range("I40").Value = 2485.2307
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 2485.2307
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -2349.2307
$ws.Range("N40").Value = -3022

$ws.Range("H55").Value = 209.63637
$ws.Range("I55").Value = 172.5
$ws.Range("J55").Value = 308.66666
$ws.Range("K55").Value = 172.5
$ws.Range("L55").Value = 308.66666
$ws.Range("M55").Value = 0.5
$ws.Range("N55").Value = -654.66666

$ws.Range("H122").Value = 2694.2646
$ws.Range("I122").Value = 2516.6667
$ws.Range("J122").Value = 3379.2856
$ws.Range("K122").Value = 7550.000100000001
$ws.Range("L122").Value = 10137.8568
$ws.Range("M122").Value = -5100.000100000001
$ws.Range("N122").Value = -15037.8568
